# BoM.xlsx update - "Se actualizan los costos de los componentes"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel PasteSpecial / alignment constants used below
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Row 2: RP2040 - price/unit update (45.69 -> 114, qty 1 -> 2)
# ---------------------------------------------------------------------
$ws.Range("D2").Value = 114
$ws.Range("E2").Value = 2

# ---------------------------------------------------------------------
# Row 4: BMP-180 - price update (49 -> 17)
# ---------------------------------------------------------------------
$ws.Range("D4").Value = 17

# ---------------------------------------------------------------------
# Row 7: used to be an (almost) empty placeholder row - becomes the new
# "Air quality" / MQ-7 gas sensor component.
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "Air quality"
$ws.Range("B7").Value = "MQ-7"

$ws.Range("D7").Value = 46
$ws.Range("D3").Copy()
$ws.Range("D7").PasteSpecial($xlPasteFormats)

$ws.Range("E7").Value = 1

$ws.Range("G7").Value = "UNIT Electronics"
$ws.Range("G3").Copy()
$ws.Range("G7").PasteSpecial($xlPasteFormats)

$ws.Hyperlinks.Add($ws.Range("H7"), "https://uelectronics.com/producto/sensor-de-gas-mq-7/", "", "", "UNIT - MQ7")
$ws.Hyperlinks.Add($ws.Range("I7"), "https://uelectronics.com/producto/sensor-de-gas-mq-7/", "", "", "MQ-7")

# ---------------------------------------------------------------------
# Row 10: LD1117AS33TR (Voltage regulator 3.3V) - price update
# ---------------------------------------------------------------------
$ws.Range("D10").Value = 23

# ---------------------------------------------------------------------
# Row 11: brand new "Headers" / 40 Pins 2mm component row
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "Headers"

$ws.Range("B11").Value = "40 Pins 2mm"
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial($xlPasteFormats)

$ws.Range("D11").Value = 6
$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial($xlPasteFormats)

$ws.Range("E11").Value = 2

$ws.Range("F10").Copy()
$ws.Range("F11").PasteSpecial($xlPasteFormats)
$ws.Range("F11").Formula = "=D11*E11"

$ws.Range("G11").Value = "UNIT Electronics"
$ws.Range("G10").Copy()
$ws.Range("G11").PasteSpecial($xlPasteFormats)

$ws.Hyperlinks.Add($ws.Range("H11"), "https://uelectronics.com/producto/header-40-pines-2-54mm/", "", "", "UNIT - Headers")
$ws.Hyperlinks.Add($ws.Range("I11"), "https://uelectronics.com/producto/header-40-pines-2-54mm/", "", "", "Headers")

# ---------------------------------------------------------------------
# Reflect the last-used selection in the sheet view
# ---------------------------------------------------------------------
$ws.Range("I11").Select()
